$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 24 (old) moves to row 23, with its Remarks text changed ---
# Clear the old row 24 entirely first (it will be rewritten as row 23).
$ws.Range("A24:I24").Clear()

# NOTE on write order: this engine's shared-string table is append-on-first-use,
# so new string cells must be written in the exact order their text first
# appears in the target file, regardless of row order, to land on the right
# <sst> index. The old row 24 (now row 23) "Remarks" cell reused its original
# slot (index 34) once the only reference to the old text was cleared above,
# so that text must be (re)written first to reclaim it, before any brand-new
# strings get appended after it.
$ws.Range("I23").Value = "Stop halfway"                        # -> reclaims index 34

$ws.Range("I25").Value = "Implemented multicore testing"      # -> new index 35
$ws.Range("I27").Value = "Same, but with ucto"                # -> new index 36
$ws.Range("I32").Value = "Wessel1, Wessel1-lex, Sonar1, Sonar1-lex"  # -> new index 37
$ws.Range("I30").Value = "Wessel1, Sonar1, Wessel1-lex, Sonar1-lex"  # -> new index 38
$ws.Range("I33").Value = "Wessel1, Wessel1-lex"                # -> new index 39
$ws.Range("I34").Value = "Sonar1, Sonar1-lex"                  # -> new index 40

# Reuses of the strings just created (order among these doesn't matter anymore).
$ws.Range("I26").Value = "Implemented multicore testing"
$ws.Range("I28").Value = "Same, but with ucto"
$ws.Range("I31").Value = "Wessel1, Sonar1, Wessel1-lex, Sonar1-lex"

$ws.Range("A23").Value = "Wessel1"
$ws.Range("B23").Value = 0.1
$ws.Range("B23").NumberFormat = "0%"
$ws.Range("C23").Value = 3
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = "w"
$ws.Range("F23").Value = 19
$ws.Range("G23").Value = 25

# --- New experiment rows ---
# Row 25
$ws.Range("A25").Value = "Wessel1"
$ws.Range("B25").Value = "Standardtest2"
$ws.Range("C25").Value = 3
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = "w"
$ws.Range("F25").Value = 9
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = 77

# Row 26
$ws.Range("A26").Value = "Wessel1"
$ws.Range("B26").Value = 0.1
$ws.Range("B26").NumberFormat = "0%"
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = "w"
$ws.Range("F26").Value = 24
$ws.Range("G26").Value = 29
$ws.Range("H26").Value = 11579

# Row 27
$ws.Range("A27").Value = "Wessel1"
$ws.Range("B27").Value = "Standardtest2"
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 3
$ws.Range("E27").Value = "w"
$ws.Range("F27").Value = 9
$ws.Range("G27").Value = 10
$ws.Range("H27").Value = 73

# Row 28
$ws.Range("A28").Value = "Wessel1"
$ws.Range("B28").Value = 0.1
$ws.Range("B28").NumberFormat = "0%"
$ws.Range("C28").Value = 3
$ws.Range("D28").Value = 3
$ws.Range("E28").Value = "w"
$ws.Range("F28").Value = 25
$ws.Range("G28").Value = 29
$ws.Range("H28").Value = 13760

# Row 29 intentionally left blank (matches source gap)

# Row 30
$ws.Range("A30").Value = "Wessel1"
$ws.Range("B30").Value = "Standardtest2"
$ws.Range("C30").Value = 3
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = "w"
$ws.Range("F30").Value = 9
$ws.Range("G30").Value = 12
$ws.Range("H30").Value = 146

# Row 31
$ws.Range("A31").Value = "Wessel1"
$ws.Range("B31").Value = 0.1
$ws.Range("B31").NumberFormat = "0%"
$ws.Range("C31").Value = 3
$ws.Range("D31").Value = 3
$ws.Range("E31").Value = "w"
$ws.Range("F31").Value = 23
$ws.Range("G31").Value = 30
$ws.Range("H31").Value = 22508

# Row 32
$ws.Range("A32").Value = "Wessel1"
$ws.Range("B32").Value = 0.1
$ws.Range("B32").NumberFormat = "0%"
$ws.Range("C32").Value = 3
$ws.Range("D32").Value = 3
$ws.Range("E32").Value = "w"
$ws.Range("F32").Value = 23
$ws.Range("G32").Value = 30
$ws.Range("H32").Value = 22292

# Row 33
$ws.Range("A33").Value = "Wessel1"
$ws.Range("B33").Value = 0.1
$ws.Range("B33").NumberFormat = "0%"
$ws.Range("C33").Value = 3
$ws.Range("D33").Value = 3
$ws.Range("E33").Value = "w"
$ws.Range("F33").Value = 23
$ws.Range("G33").Value = 29
$ws.Range("H33").Value = 20800

# Row 34
$ws.Range("A34").Value = "Wessel1"
$ws.Range("B34").Value = 0.1
$ws.Range("B34").NumberFormat = "0%"
$ws.Range("C34").Value = 3
$ws.Range("D34").Value = 3
$ws.Range("E34").Value = "w"
$ws.Range("F34").Value = 8
$ws.Range("G34").Value = 10
$ws.Range("H34").Value = 31446

# --- Update the view: scroll position + selection ---
$ws.Application.ActiveWindow.ScrollRow = 16
$ws.Range("G34").Select()
